$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking": Right (B11) 5 -> 4, Wrong (C11) -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right (B12) 140 -> 112, and Max label (E12) text updated accordingly
$ws.Range("B12").Value = 112
$ws.Range("E12").Value = "112 / 112"
